$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Trening" header in column F, reusing the same style as the
# other header cells (copy format from E1, which carries the bold/border
# header style already present in the workbook).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "Trening"

# New data rows (A: timestamp serial, B: Seconds, C: Velocity,
# D: Acceleration_SMA, E: Velocity_Bin, F: Trening)
$data = @(
    @(45684.59215335648, 562, 10.65, 1.596183095659528, "10-15", "Duża Gra"),
    @(45684.59280266204, 618.1, 11.27, 2.449963876179286, "10-15", "Duża Gra"),
    @(45684.5941244213, 732.3, 10.9, 2.008010251181468, "10-15", "Duża Gra"),
    @(45684.59280034722, 617.9, 9.27, 2.296086515699114, "5-10", "Duża Gra"),
    @(45684.59403761574, 724.8, 9.9, 1.596917390823366, "5-10", "Duża Gra"),
    @(45684.59412210648, 732.1, 9.34, 1.946478349821909, "5-10", "Duża Gra"),
    @(45684.5991417824, 1165.8, 11.61, 3.213301862989153, "10-15", "Mała Gra"),
    @(45684.60245775463, 1452.3, 14.21, 3.690788303102764, "10-15", "Mała Gra"),
    @(45684.60405729167, 1590.5, 10.97, 3.230537108012609, "10-15", "Mała Gra"),
    @(45684.59671701389, 956.3, 8.92, 2.957239730017526, "5-10", "Mała Gra"),
    @(45684.59913946759, 1165.6, 9.140000000000001, 3.096880299704416, "5-10", "Mała Gra"),
    @(45684.60405613426, 1590.4, 9.789999999999999, 3.15067059653146, "5-10", "Mała Gra")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

# Write the timestamp as a real date serial number. Applying the number
# format twice on the first cell (once with a lowercase format code, then
# the uppercase one actually used) mirrors how the style table ends up
# with both format codes registered while only the uppercase one is
# referenced by any cell style.
$ws.Cells.Item(2, 1).Value = $data[0][0]
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 1; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
}
$ws.Range("A3:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
